$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.072.59'
$ws.Range("E2").Value = '  -1.33%  '

# Row 3
$ws.Range("D3").Value = '2.687.40'
$ws.Range("E3").Value = '  -1.90%  '

# Row 4
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").Value = "'556.72"
$ws.Range("E5").Value = '  -1.78%  '

# Row 6
$ws.Range("D6").Value = "'159.18"
$ws.Range("E6").Value = '  -0.44%  '

# Row 7
$ws.Range("E7").Value = '  -0.06%  '

# Row 8
$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = '  -0.61%  '

# Row 9
$ws.Range("E9").Value = '  -2.87%  '

# Row 10
$ws.Range("E10").Value = '  -1.98%  '

# Row 11
$ws.Range("D11").Value = "'0.371"
$ws.Range("E11").Value = '  -3.32%  '

# Row 12
$ws.Range("D12").Value = "'5.39"
$ws.Range("E12").Value = '  -6.36%  '

# Row 13
$ws.Range("D13").Value = '3.162.75'
$ws.Range("E13").Value = '  -2.03%  '

# Row 14
$ws.Range("D14").Value = "'26.57"
$ws.Range("E14").Value = '  -1.60%  '

# Row 15
$ws.Range("D15").Value = '62.966.46'
$ws.Range("E15").Value = '  -1.32%  '

# Row 16
$ws.Range("D16").Value = "'0.0000148"
$ws.Range("E16").Value = '  -1.31%  '

# Row 17
$ws.Range("D17").Value = '2.687.15'
$ws.Range("E17").Value = '  -2.09%  '

# Row 18
$ws.Range("D18").Value = "'12.03"
$ws.Range("E18").Value = '  -1.38%  '

# Row 19
$ws.Range("D19").Value = "'4.64"
$ws.Range("E19").Value = '  -3.40%  '

# Row 20
$ws.Range("D20").Value = "'346.59"
$ws.Range("E20").Value = '  -2.08%  '

# Row 21
$ws.Range("D21").Value = "'6.31"
$ws.Range("E21").Value = '  -4.52%  '

# Row 22
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = '  +0.36%  '

# Row 23
$ws.Range("D23").Value = "'0.514"
$ws.Range("E23").Value = '  -1.61%  '

# Row 24
$ws.Range("D24").Value = "'63.41"
$ws.Range("E24").Value = '  -1.41%  '

# Row 25
$ws.Range("E25").Value = '  -0.91%  '

# Row 26
$ws.Range("E26").Value = '  +0.17%  '

# Row 27
$ws.Range("D27").Value = "'8.25"
$ws.Range("E27").Value = '  -2.04%  '

# Row 28
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").Value = "'1.43"
$ws.Range("E28").Value = '  +7.69%  '

# Row 29
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0867'
$ws.Range("E29").Value = '  -4.66%  '

# Row 30
$ws.Range("E30").Value = '  +0.86%  '

# Row 31
$ws.Range("E31").Value = '  -0.47%  '

# Row 32
$ws.Range("D32").Value = "'165.43"
$ws.Range("E32").Value = '  +0.90%  '

# Row 33
$ws.Range("D33").Value = "'4.96"
$ws.Range("E33").Value = '  +1.25%  '

# Row 34
$ws.Range("D34").Value = "'1.50"
$ws.Range("E34").Value = '  +1.56%  '

# Row 35
$ws.Range("E35").Value = '  -0.01%  '

# Row 36
$ws.Range("D36").Value = "'19.55"
$ws.Range("E36").Value = '  -2.71%  '

# Row 37
$ws.Range("E37").Value = '  -0.22%  '

# Row 38
$ws.Range("D38").Value = "'358.13"
$ws.Range("E38").Value = '  +2.80%  '

# Row 39
$ws.Range("D39").Value = "'6.42"
$ws.Range("E39").Value = '  +1.27%  '

# Row 40
$ws.Range("D40").Value = "'0.962"
$ws.Range("E40").Value = '  -2.40%  '

# Row 41
$ws.Range("E41").Value = '  -1.87%  '

# Row 42
$ws.Range("D42").Value = "'38.52"
$ws.Range("E42").Value = '  -0.26%  '

# Row 43
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = "'21.07"
$ws.Range("E43").Value = '  -3.73%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'20.48"
$ws.Range("E44").Value = '  -2.79%  '

# Row 45
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").Value = "'0.0564"
$ws.Range("E45").Value = '  -3.36%  '

# Row 46
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = "'0.620"
$ws.Range("E46").Value = '  -0.61%  '

# Row 47
$ws.Range("E47").Value = '  -0.09%  '

# Row 48
$ws.Range("E48").Value = '  +0.08%  '

# Row 49
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = "'0.0244"
$ws.Range("E49").Value = '  -2.54%  '

# Row 50
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = "'0.0975"
$ws.Range("E50").Value = '  -2.83%  '

# Row 51
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = "'129.80"
$ws.Range("E51").Value = '  -3.44%  '
